$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.266.25'
$ws.Range("E2").Value = '  +6.59%  '
$ws.Range("D3").Value = '3.671.72'
$ws.Range("E3").Value = '  +18.58%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '606.91'
$ws.Range("E5").Value = '  +5.34%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '182.03'
$ws.Range("E6").Value = '  +2.81%  '
$ws.Range("D7").Value = '3.668.61'
$ws.Range("E7").Value = '  +18.55%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("E9").Value = '  +5.57%  '
$ws.Range("E10").Value = '  +8.23%  '
$ws.Range("E11").Value = '  +4.86%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.499'
$ws.Range("E12").Value = '  +7.69%  '
$ws.Range("E13").Value = '  +13.42%  '
$ws.Range("E14").Value = '  +6.37%  '
$ws.Range("D15").Value = '4.287.11'
$ws.Range("E15").Value = '  +18.52%  '
$ws.Range("D16").Value = '71.255.05'
$ws.Range("E16").Value = '  +6.54%  '
$ws.Range("D17").Value = '3.675.27'
$ws.Range("E17").Value = '  +18.38%  '
$ws.Range("E18").Value = '  +1.19%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.52'
$ws.Range("E19").Value = '  +8.20%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.02'
$ws.Range("E20").Value = '  +1.03%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '519.30'
$ws.Range("E21").Value = '  +7.92%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.26'
$ws.Range("E22").Value = '  +20.38%  '
$ws.Range("E23").Value = '  +8.62%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '88.29'
$ws.Range("E24").Value = '  +5.83%  '
$ws.Range("E25").Value = '  +11.72%  '
$ws.Range("E26").Value = '  +7.68%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.99'
$ws.Range("E27").Value = '  +8.17%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("E29").Value = '  +12.23%  '
$ws.Range("E30").Value = '  +1.54%  '
$ws.Range("E31").Value = '  +8.37%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0000111'
$ws.Range("E32").Value = '  +18.65%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.75'
$ws.Range("E33").Value = '  +13.60%  '
$ws.Range("E34").Value = '  +5.11%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  -0.05%  '
$ws.Range("E36").Value = '  +10.22%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.01'
$ws.Range("E37").Value = '  +8.00%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.347'
$ws.Range("E38").Value = '  +12.55%  '
$ws.Range("E39").Value = '  +12.00%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '51.53'
$ws.Range("E40").Value = '  +5.16%  '
$ws.Range("E41").Value = '  +6.16%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '45.14'
$ws.Range("E42").Value = '  -6.75%  '
$ws.Range("B43").Value = 'Cosmos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.84'
$ws.Range("E43").Value = '  +7.18%  '
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '3.131.83'
$ws.Range("E44").Value = '  +12.37%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '417.77'
$ws.Range("E45").Value = '  +12.68%  '
$ws.Range("E46").Value = '  +4.94%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '28.56'
$ws.Range("E47").Value = '  +14.97%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0369'
$ws.Range("E48").Value = '  +7.68%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '139.23'
$ws.Range("E49").Value = '  +2.82%  '
$ws.Range("E50").Value = '  -0.01%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.47'
$ws.Range("E51").Value = '  +11.79%  '
